# Update "想去人数" (want-to-go count) values in column F across the
# "展览", "演出", "本地生活" and "全部类型" sheets, matching the gh-pages
# regeneration captured at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value  = 74
$ws1.Cells.Item(3, 6).Value  = 197
$ws1.Cells.Item(4, 6).Value  = 58
$ws1.Cells.Item(5, 6).Value  = 1634
$ws1.Cells.Item(6, 6).Value  = 3243
$ws1.Cells.Item(7, 6).Value  = 785
$ws1.Cells.Item(8, 6).Value  = 2000
$ws1.Cells.Item(9, 6).Value  = 1922
$ws1.Cells.Item(10, 6).Value = 982
$ws1.Cells.Item(12, 6).Value = 15
$ws1.Cells.Item(13, 6).Value = 1589
$ws1.Cells.Item(14, 6).Value = 339
$ws1.Cells.Item(16, 6).Value = 63
$ws1.Cells.Item(17, 6).Value = 30
$ws1.Cells.Item(18, 6).Value = 1404
$ws1.Cells.Item(19, 6).Value = 500
$ws1.Cells.Item(20, 6).Value = 604
$ws1.Cells.Item(21, 6).Value = 302
$ws1.Cells.Item(22, 6).Value = 10534
$ws1.Cells.Item(23, 6).Value = 9714
$ws1.Cells.Item(24, 6).Value = 839
$ws1.Cells.Item(25, 6).Value = 644
$ws1.Cells.Item(26, 6).Value = 1813
$ws1.Cells.Item(27, 6).Value = 135
$ws1.Cells.Item(28, 6).Value = 393

# Sheet "演出" (sheetId 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(4, 6).Value = 35
$ws2.Cells.Item(5, 6).Value = 114

# Sheet "本地生活" (sheetId 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 58

# Sheet "全部类型" (sheetId 4) - combined listing of all the above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value  = 74
$ws4.Cells.Item(3, 6).Value  = 58
$ws4.Cells.Item(4, 6).Value  = 197
$ws4.Cells.Item(6, 6).Value  = 58
$ws4.Cells.Item(7, 6).Value  = 1634
$ws4.Cells.Item(8, 6).Value  = 3243
$ws4.Cells.Item(9, 6).Value  = 785
$ws4.Cells.Item(10, 6).Value = 2000
$ws4.Cells.Item(11, 6).Value = 1923
$ws4.Cells.Item(12, 6).Value = 982
$ws4.Cells.Item(14, 6).Value = 15
$ws4.Cells.Item(15, 6).Value = 1589
$ws4.Cells.Item(16, 6).Value = 339
$ws4.Cells.Item(19, 6).Value = 63
$ws4.Cells.Item(20, 6).Value = 35
$ws4.Cells.Item(21, 6).Value = 30
$ws4.Cells.Item(22, 6).Value = 1404
$ws4.Cells.Item(23, 6).Value = 500
$ws4.Cells.Item(24, 6).Value = 604
$ws4.Cells.Item(25, 6).Value = 302
$ws4.Cells.Item(26, 6).Value = 10534
$ws4.Cells.Item(27, 6).Value = 9714
$ws4.Cells.Item(28, 6).Value = 839
$ws4.Cells.Item(29, 6).Value = 644
$ws4.Cells.Item(30, 6).Value = 1813
$ws4.Cells.Item(31, 6).Value = 114
$ws4.Cells.Item(33, 6).Value = 135
$ws4.Cells.Item(34, 6).Value = 393
